# actualiza dashboard y graficos
# Strip the survey-stage suffix (e.g. "_ELE", "_ENI", "_ENIA", "_EOC") from the
# "responsible" column (C) values on rows 2-23 of the active sheet, e.g.
# "CamiloR_ELE" -> "CamiloR", "CatalinaJ_ENI" -> "CatalinaJ", etc.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 23; $r++) {
    $cell = $ws.Range("C$r")
    $old = $cell.Text
    if ($old -match '^(.*)_[A-Za-z]+$') {
        $cell.Value = $matches[1]
    }
}
